# Dream Team Project 1
#
# Split the "Approximately 30 million ... shoots(Statista) ." bullet on the
# "Facts" slide into three runs so the middle portion reads
# "shoots in 2017(" instead of "shoots(", i.e. the bullet becomes:
#   "Approximately 30 million participating in firearms target "
#   "shoots in 2017("
#   "Statista) ."

$p = $ppt.ActivePresentation

$oldFull = "Approximately 30 million participating in firearms target shoots(Statista) ."
$searchKey = "Approximately 30 million"

$targetShape = $null
$targetSlide = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $txt = $tf.TextRange.Text
                if ($txt -like "*$searchKey*") {
                    $targetShape = $shp
                    $targetSlide = $slide
                }
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Locate the exact original sentence so we know precisely where it starts.
$found = $tr.Find($oldFull)
$start = $found.Start

# Within that run, "shoots(" (7 characters) sits right after the 58-character
# lead-in "Approximately 30 million participating in firearms target ".
# Replacing just those 7 characters with "shoots in 2017(" naturally splits
# the original single run into three runs - the unchanged lead-in, the
# replaced/inserted piece, and the unchanged "Statista) ." tail - while each
# keeps the existing Cambria font formatting.
$leadIn = "Approximately 30 million participating in firearms target "
$middleOld = "shoots("
$middleNew = "shoots in 2017("

$midRange = $tr.Characters($start + $leadIn.Length, $middleOld.Length)
$midRange.Text = $middleNew
